# Rename the worksheet from the German default "Tabelle1" to the
# English default "Sheet1" (the only content-level change in the target
# revision; everything else in the recorded diff -- the revisionPtr GUID,
# the bookViews window position, and the cellXfs ordering in styles.xml --
# is incidental re-save noise with no visible effect on any cell).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Sheet1"
